# Add payment 79174445 (Cash) 2025-08-18T17:46:33 as a new row (48) at the
# bottom of the payments sheet, mirroring the existing rows for this phone
# number (e.g. row 45: original_amount == final_amount, no discount/points).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 47 stored the phone number as text; normalize it to a
# real number so it is consistent with the rest of column A.
$ws.Cells.Item(47, 1).Value = 79174445

$row = 48

$ws.Cells.Item($row, 1).Value = 79174445            # A: phone
$ws.Cells.Item($row, 2).Value = ""                  # B: amount (blank)
$ws.Cells.Item($row, 3).Value = "Cash"               # C: method
$ws.Cells.Item($row, 4).Value = "2025-08-18T17:46:33" # D: timestamp
$ws.Cells.Item($row, 5).Value = 200                  # E: original_amount
$ws.Cells.Item($row, 6).Value = ""                  # F: discount_applied (blank)
$ws.Cells.Item($row, 7).Value = 200                  # G: final_amount
$ws.Cells.Item($row, 8).Value = 0                    # H: birthday_discount
$ws.Cells.Item($row, 9).Value = 0                    # I: points_redeemed
